$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPriceValue (optional); E = newVolumeValue (optional) }
$updates = @{
    2  = @{ D = "36.605.56";    E = "  -0.42%  " }
    3  = @{ D = "1.967.82";     E = "  +0.24%  " }
    4  = @{ E = "  +0.03%  " }
    5  = @{ E = "  -0.08%  " }
    6  = @{ E = "  +1.82%  " }
    7  = @{ D = "59.97";        E = "  +2.52%  " }
    8  = @{ E = "  +0.01%  " }
    9  = @{ E = "  +2.03%  " }
    10 = @{ E = "  -2.81%  " }
    12 = @{ D = "14.24";        E = "  +3.87%  " }
    13 = @{ D = "0.840";        E = "  +2.35%  " }
    14 = @{ D = "2.257.67";     E = "  +0.52%  " }
    15 = @{ D = "21.59";        E = "  -2.53%  " }
    16 = @{ E = "  +0.26%  " }
    17 = @{ D = "1.972.49";     E = "  +0.60%  " }
    18 = @{ D = "36.530.23";    E = "  -0.49%  " }
    19 = @{ D = "69.72";        E = "  -0.02%  " }
    20 = @{ D = ([string]"0.0" + [string][char]0x2083 + [string]"0854"); E = "  -1.13%  " }
    21 = @{ D = "229.31";       E = "  +0.51%  " }
    22 = @{ D = "5.07";         E = "  -0.50%  " }
    23 = @{ E = "  +0.05%  " }
    24 = @{ E = "  +1.50%  " }
    25 = @{ E = "  +1.17%  " }
    26 = @{ D = "0.146";        E = "  +6.08%  " }
    27 = @{ D = "9.12";         E = "  -1.94%  " }
    28 = @{ D = "162.37" }
    29 = @{ E = "  +0.13%  " }
    30 = @{ E = "  +20.35%  " }
    31 = @{ E = "  +2.17%  " }
    32 = @{ E = "  +2.63%  " }
    33 = @{ E = "  -0.71%  " }
    34 = @{ E = "  +6.87%  " }
    35 = @{ E = "  +3.91%  " }
    36 = @{ E = "  +0.09%  " }
    37 = @{ D = "3.35";         E = "  -2.76%  " }
    38 = @{ E = "  +0.24%  " }
    39 = @{ D = "5.40";         E = "  -13.67%  " }
    40 = @{ D = "0.0962";       E = "  -3.97%  " }
    41 = @{ E = "  +0.54%  " }
    42 = @{ E = "  -0.13%  " }
    43 = @{ D = "0.0210";       E = "  -1.37%  " }
    44 = @{ D = "15.87";        E = "  -1.07%  " }
    45 = @{ D = "1.366.27";     E = "  +1.49%  " }
    46 = @{ D = "88.95";        E = "  +1.73%  " }
    47 = @{ E = "  -1.59%  " }
    48 = @{ E = "  +0.84%  " }
    49 = @{ D = "2.81";         E = "  -0.91%  " }
    50 = @{ D = "45.83";        E = "  +5.24%  " }
    51 = @{ D = "2.151.43";     E = "  +0.63%  " }
}

# Rows whose new Price text looks like a plain number to Excel's type
# sniffer (e.g. "59.97") need NumberFormat forced to Text first, or the
# COM layer silently stores them as a numeric value instead of a string -
# which would not match the original inlineStr/text cell content.
# Rows whose text has two dots (e.g. "36.605.56") or other non-numeric
# characters are already safe and are left with the default style so we
# don't disturb formatting that doesn't need to change.
$needsTextFormat = @(7, 12, 13, 15, 19, 21, 22, 26, 27, 28, 37, 39, 40, 43, 44, 46, 49, 50)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        if ($needsTextFormat -contains [int]$row) {
            # Force the cell to Text so the numeric-looking string isn't
            # silently reinterpreted as a number, then restore the
            # original cell style so no stray formatting change lingers.
            $origStyle = $cell.Style
            $cell.NumberFormat = "@"
            $cell.Value = $vals["D"]
            $cell.Style = $origStyle
        } else {
            $cell.Value = $vals["D"]
        }
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
